$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 ("Diseases (patient-stated)"), shifting rows 2-12 down to 3-13.
$ws.Rows("2:2").Insert()

# Delete the row that is now pushed past the bottom of the table (old row 12 -> now row 13,
# "Symptom – Skin/Hair/Nails"), since the table keeps only 12 rows total (1 header + 11 data rows... )
$ws.Rows("13:13").Delete()

# New row 2 values: "Abnormal test result"
$ws.Range("A2").Value = "Abnormal test result"
$ws.Range("B2").Value = 1.5
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1.6
$ws.Range("E2").Value = 2.5
$ws.Range("F2").Value = 2.8
$ws.Range("G2").Value = 2.4

# Copy the style of the label cell from a neighboring row so it matches the others (bold/border/center).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Reset the numeric cells in the new row back to the plain (default) style, since the row-insert
# operation can carry formatting down from the header row above.
$ws.Range("B2:G2").ClearFormats()

# Row 3: Diseases (patient-stated)
$ws.Range("B3").Value = 6.2
$ws.Range("C3").Value = 5.6
$ws.Range("D3").Value = 4.9
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 6.9
$ws.Range("G3").Value = 6.3

# Row 4: Injuries & adverse effects
$ws.Range("B4").Value = 12.1
$ws.Range("C4").Value = 11.1
$ws.Range("D4").Value = 7.9
$ws.Range("E4").Value = 9.300000000000001
$ws.Range("F4").Value = 7.2
$ws.Range("G4").Value = 5

# Row 5: Other
$ws.Range("B5").Value = 9.199999999999999
$ws.Range("C5").Value = 11.1
$ws.Range("D5").Value = 6.5
$ws.Range("E5").Value = 8.4
$ws.Range("F5").Value = 6.4
$ws.Range("G5").Value = 6.4

# Row 6: Symptom – Circulatory
$ws.Range("B6").Value = 8.6
$ws.Range("C6").Value = 27.8
$ws.Range("D6").Value = 7.5
$ws.Range("E6").Value = 10.1
$ws.Range("F6").Value = 10.8
$ws.Range("G6").Value = 9.4

# Row 7: Symptom – Digestive
$ws.Range("B7").Value = 12.8
$ws.Range("C7").Value = 5.6
$ws.Range("D7").Value = 9.699999999999999
$ws.Range("E7").Value = 13.6
$ws.Range("F7").Value = 13.6
$ws.Range("G7").Value = 16.8

# Row 8: Symptom – Eye/Ear
$ws.Range("B8").Value = 5.3
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 4.8
$ws.Range("E8").Value = 6.7
$ws.Range("F8").Value = 7.6
$ws.Range("G8").Value = 6.6

# Row 9: Symptom – General
$ws.Range("B9").Value = 6.6
$ws.Range("C9").Value = 5.6
$ws.Range("D9").Value = 4.8
$ws.Range("E9").Value = 6.6
$ws.Range("F9").Value = 5.5
$ws.Range("G9").Value = 5.9

# Row 10: Symptom – Genitourinary
$ws.Range("B10").Value = 9.800000000000001
$ws.Range("C10").Value = 5.6
$ws.Range("D10").Value = 11.3
$ws.Range("E10").Value = 16.3
$ws.Range("F10").Value = 15.2
$ws.Range("G10").Value = 14.1

# Row 11: Symptom – Nervous
$ws.Range("B11").Value = 24.9
$ws.Range("C11").Value = 16.7
$ws.Range("D11").Value = 38.5
$ws.Range("E11").Value = 14.6
$ws.Range("F11").Value = 21.8
$ws.Range("G11").Value = 24.9

# Row 12: Symptom – Respiratory
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 11.1
$ws.Range("D12").Value = 2.4
$ws.Range("E12").Value = 2.9
$ws.Range("F12").Value = 2.1
$ws.Range("G12").Value = 2.1
